$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.009821333333333333
$ws.Range("H2").Value = 0.029464
$ws.Range("I2").Value = 0.06297798848338983
$ws.Range("J2").Value = 0.06297798848338984
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2615913333333333
$ws.Range("N2").Value = 0.784774
$ws.Range("O2").Value = 0.08239613548481725
$ws.Range("P2").Value = 0.08239613548481727
$ws.Range("Q2").Value = 0.002569175681777778
$ws.Range("R2").Value = 0.023122581136
$ws.Range("S2").Value = 0.005189142871638648
$ws.Range("T2").Value = 0.005189142871638651

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.009821333333333333
$ws.Range("H3").Value = 0.029464
$ws.Range("I3").Value = 0.06297798848338983
$ws.Range("J3").Value = 0.06297798848338984
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.5494561706387266
$ws.Range("P3").Value = 0.5494561706387268
$ws.Range("Q3").Value = 0.01713247136533333
$ws.Range("R3").Value = 0.154192242288
$ws.Range("S3").Value = 0.0346036443866132
$ws.Range("T3").Value = 0.03460364438661322

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.009821333333333333
$ws.Range("H4").Value = 0.029464
$ws.Range("I4").Value = 0.06297798848338983
$ws.Range("J4").Value = 0.06297798848338984
$ws.Range("M4").Value = 1.168795666666667
$ws.Range("N4").Value = 3.506387
$ws.Range("O4").Value = 0.3681476938764561
$ws.Range("P4").Value = 0.3681476938764561
$ws.Range("Q4").Value = 0.01147913184088889
$ws.Range("R4").Value = 0.103312186568
$ws.Range("S4").Value = 0.02318520122513797
$ws.Range("T4").Value = 0.02318520122513798

# Row 5
$ws.Range("I5").Value = 0.3247949111459754
$ws.Range("J5").Value = 0.3247949111459754
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2615913333333333
$ws.Range("N5").Value = 0.784774
$ws.Range("O5").Value = 0.08239613548481725
$ws.Range("P5").Value = 0.08239613548481727
$ws.Range("Q5").Value = 0.01324994982177778
$ws.Range("R5").Value = 0.119249548396
$ws.Range("S5").Value = 0.02676184550356297
$ws.Range("T5").Value = 0.02676184550356298

# Row 6
$ws.Range("I6").Value = 0.3247949111459754
$ws.Range("J6").Value = 0.3247949111459754
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.5494561706387266
$ws.Range("P6").Value = 0.5494561706387268
$ws.Range("R6").Value = 0.7952120548680002
$ws.Range("S6").Value = 0.1784605681212131
$ws.Range("T6").Value = 0.1784605681212132

# Row 7
$ws.Range("I7").Value = 0.3247949111459754
$ws.Range("J7").Value = 0.3247949111459754
$ws.Range("M7").Value = 1.168795666666667
$ws.Range("N7").Value = 3.506387
$ws.Range("O7").Value = 0.3681476938764561
$ws.Range("P7").Value = 0.3681476938764561
$ws.Range("Q7").Value = 0.05920105891088889
$ws.Range("R7").Value = 0.532809530198
$ws.Range("S7").Value = 0.1195724975211993
$ws.Range("T7").Value = 0.1195724975211993

# Row 8
$ws.Range("G8").Value = 0.09547600000000001
$ws.Range("H8").Value = 0.286428
$ws.Range("I8").Value = 0.6122271003706348
$ws.Range("J8").Value = 0.6122271003706349
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2615913333333333
$ws.Range("N8").Value = 0.784774
$ws.Range("O8").Value = 0.08239613548481725
$ws.Range("P8").Value = 0.08239613548481727
$ws.Range("Q8").Value = 0.02497569414133334
$ws.Range("R8").Value = 0.224781247272
$ws.Range("S8").Value = 0.05044514710961563
$ws.Range("T8").Value = 0.05044514710961565

# Row 9
$ws.Range("G9").Value = 0.09547600000000001
$ws.Range("H9").Value = 0.286428
$ws.Range("I9").Value = 0.6122271003706348
$ws.Range("J9").Value = 0.6122271003706349
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.5494561706387266
$ws.Range("P9").Value = 0.5494561706387268
$ws.Range("Q9").Value = 0.166549671064
$ws.Range("R9").Value = 1.498947039576
$ws.Range("S9").Value = 0.3363919581309003
$ws.Range("T9").Value = 0.3363919581309004

# Row 10
$ws.Range("G10").Value = 0.09547600000000001
$ws.Range("H10").Value = 0.286428
$ws.Range("I10").Value = 0.6122271003706348
$ws.Range("J10").Value = 0.6122271003706349
$ws.Range("M10").Value = 1.168795666666667
$ws.Range("N10").Value = 3.506387
$ws.Range("O10").Value = 0.3681476938764561
$ws.Range("P10").Value = 0.3681476938764561
$ws.Range("Q10").Value = 0.1115919350706667
$ws.Range("R10").Value = 1.004327415636
$ws.Range("S10").Value = 0.2253899951301188
$ws.Range("T10").Value = 0.2253899951301189
